$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 corresponds to model "meta-llama/llama-4-scout-17b-16e-instruct"
# Update Current_Ct_Day (I5), Current_Pct_Ct (J5), Current_Ct_Tokens (K5), Current_Pct_Tokens (L5)
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 0.002
$ws.Range("K5").Value = 948
$ws.Range("L5").Value = 0.001896
